$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D; existing data in D:K shifts right to E:L
$ws.Columns("D:D").Insert()

# Copy formatting for the new column D from column E (the just-shifted former column D)
# so the new cells inherit the correct number formats/styles, per contiguous data block
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new column D with the latest (FY2018) financial figures
$ws.Range("D7").Value2 = 43465
$ws.Range("D8").Value2 = 10231000
$ws.Range("D9").Value2 = 2061000
$ws.Range("D10").Value2 = 8170000
$ws.Range("D12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = -146000
$ws.Range("D15").Value2 = 1145000
$ws.Range("D17").Value2 = 9086000
$ws.Range("D18").Value2 = 1145000
$ws.Range("D20").Value2 = 298000
$ws.Range("D21").Value2 = 2588000
$ws.Range("D22").Value2 = 576000
$ws.Range("D23").Value2 = 867000
$ws.Range("D24").Value2 = -10000
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 877000
$ws.Range("D27").Value2 = 775000
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = 0
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = -298000
$ws.Range("D33").Value2 = 775000
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 775000
$ws.Range("D38").Value2 = 43465
$ws.Range("D41").Value2 = 4000
$ws.Range("D42").Value2 = 0
$ws.Range("D43").Value2 = 1247000
$ws.Range("D44").Value2 = 0
$ws.Range("D45").Value2 = 347000
$ws.Range("D46").Value2 = 1598000
$ws.Range("D47").Value2 = "NA"
$ws.Range("D48").Value2 = 9015000
$ws.Range("D49").Value2 = 0
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 334000
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 10947000
$ws.Range("D57").Value2 = 763000
$ws.Range("D58").Value2 = 381000
$ws.Range("D59").Value2 = 1684000
$ws.Range("D60").Value2 = 2828000
$ws.Range("D61").Value2 = 7341000
$ws.Range("D62").Value2 = 311000
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 10603000
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 1671000
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = -15660000
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = -1327000
$ws.Range("D77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("D81").Value2 = 775000
$ws.Range("D83").Value2 = 1145000
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 2000000
$ws.Range("D91").Value2 = -2267000
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = 185000
$ws.Range("D96").Value2 = -92000
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = -2186000
$ws.Range("D101").Value2 = 0
$ws.Range("D102").Value2 = -1000
